$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Basso reference: citekey year 2014 -> 2004 (content/year inside stays 2004)
$newBasso = " @article{Basso:2004,`n author = {Basso, Bruno and Cammarano, Davide and De Vita, Pasquale},`n year = {2004},`n month = {01},`n pages = {36-53},`n title = {Remotely sensed vegetation indices: theory and application for crop management},`n volume = {1},`n booktitle = {Rivista Italiana di Agrometeorologia}`n}"
$ws.Range("C7").Value2 = $newBasso

# 2. Highlight column C for rows 5-9 with the same green fill used in rows 2-4
$rng = $ws.Range("C5:C9")
$rng.Interior.Color = 5287936

# 3. Update the frozen-pane scroll position / active selection to C10
$ws.Range("C10").Select()
